# Rename the embedded logo pictures that live in the document's headers
# and footers.
#
#   * The Pearson Edexcel logo (alt text / description ends in
#     "PearsonLogo.png") is renamed from "image1.png" to "image2.png" in
#     both footers that carry it.
#   * The BTEC logo (alt text / description "BTec_Logo-Orange") is
#     renamed from "image2.jpg" to "image1.jpg" in the header that
#     carries it.
#
# Inline pictures dropped into a header/footer aren't reachable through
# Document.InlineShapes (that collection only covers the main body
# story) - they have to be found via each Section's Headers/Footers
# collection, walking every header/footer slot (primary, first page,
# even page) and renaming any picture whose description matches.

$d = $word.ActiveDocument

function Rename-LogoShapes($range) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $desc = $shape.AlternativeText
        if ($desc -like "*PearsonLogo.png") {
            $shape.Name = "image2.png"
        } elseif ($desc -eq "BTec_Logo-Orange") {
            $shape.Name = "image1.jpg"
        }
    }
}

foreach ($section in $d.Sections) {
    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            Rename-LogoShapes $header.Range
        }
    }
    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            Rename-LogoShapes $footer.Range
        }
    }
}
